$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.11639999999999
$ws.Range("E4").Value = 12.7192
$ws.Range("E6").Value = 12.0672
$ws.Range("B9").Value = 8.230200000000004
$ws.Range("E10").Value = 11.804
$ws.Range("C11").Value = -13.36599999999999
$ws.Range("E11").Value = 13.48099999999999
$ws.Range("B18").Value = 4.668400000000003
$ws.Range("B20").Value = 5.5933
$ws.Range("D21").Value = -7.554600000000002
$ws.Range("E21").Value = 13.44600000000001
